$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------
# Sheet1: "Risky" arm (row 9) now uses a lower Attr1 value (0.05
# instead of 0.34), labelled with a new shared string, and a note
# describing the varied range is added at N19.
# ---------------------------------------------------------------
$ws1.Range("D9").Value = "Risky"
$ws1.Range("F9").Value = 0.05
$ws1.Range("N19").Value = "0.05 to 0.34"

# ---------------------------------------------------------------
# Sheet2: the U-column "feasibility vs suboptimality" check now
# compares against twice the suboptimality threshold, and the
# shared formula block grows to include row 11.
# ---------------------------------------------------------------
$ws2.Range("U7").Formula = "=O7>P7*2"
$ws2.Range("U8").Formula = "=O8>P8*2"
$ws2.Range("U9").Formula = "=O9>P9*2"
$ws2.Range("U10").Formula = "=O10>P10*2"
$ws2.Range("U11").Formula = "=O11>P11*2"

# ---------------------------------------------------------------
# Sheet3: Attr2/Attr4 values for the best arm (row 7) are bumped
# up, and the U-column feasibility check mirrors the Sheet2 change
# (threshold doubled, shared formula extended through row 11).
# ---------------------------------------------------------------
$ws3.Range("G7").Value = 0.9
$ws3.Range("I7").Value = 0.8
$ws3.Range("U7").Formula = "=O7>P7*2"
$ws3.Range("U8").Formula = "=O8>P8*2"
$ws3.Range("U9").Formula = "=O9>P9*2"
$ws3.Range("U10").Formula = "=O10>P10*2"
$ws3.Range("U11").Formula = "=O11>P11*2"

# ---------------------------------------------------------------
# Restore the active-cell selections recorded in each sheet view.
# Sheet3 is selected last so it remains the active tab, matching
# the workbook's saved view state.
# ---------------------------------------------------------------
$ws1.Range("N7").Select()
$ws2.Range("T17").Select()
$ws3.Range("G10").Select()
